$wb = $excel.ActiveWorkbook

# --- Sheet: VENTA MENSUAL ---
# August sales for LOZANO MOLINA TITO / PAREDES ORTIZ MARIA INES went from 0 to 2702.94
$wsVentaMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsVentaMensual.Range("F18").Value = 2702.94
$wsVentaMensual.Range("F29").Value = 4540.82

# --- Sheet: VENTAS POR GRUPO ---
# PIEDRA SINTERIZADA column for the same client/advisor pair reflects the same new sale
$wsVentasPorGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentasPorGrupo.Range("L18").Value = 2702.94
$wsVentasPorGrupo.Range("L29").Value = "2 de 27"

# --- Sheet: CUMPLIMIENTO MENSUAL ---
# PIEDRA SINTERIZADA row and TOTAL row recompute with the new sale figure
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumplimiento.Range("D15").Value = 4373.37
$wsCumplimiento.Range("E15").Value = -2735.37
$wsCumplimiento.Range("F15").Value = 2.669945054945055

$wsCumplimiento.Range("D19").Value = 4540.82
$wsCumplimiento.Range("E19").Value = 32959.18093005039
$wsCumplimiento.Range("F19").Value = 0.1210885303301751

# Note: Excel's ColumnWidth setter stores width in pixel-quantized units, so a
# nominal ColumnWidth of 21.1 is what yields the persisted width of 22 seen in
# the target workbook (mirrors Excel's own auto-fit side effect here).
$wsCumplimiento.Columns.Item(5).ColumnWidth = 21.1
